$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '69.514.02'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.693.06'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '681.50'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '160.87'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +0.25%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.147'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.44%  '
$ws.Range('E10').Value = '  -0.65%  '
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0000233'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.315.74'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '32.45'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.37%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.691.87'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.12%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '69.459.64'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('E17').Value = '  +2.63%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '16.04'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.47'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '475.48'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.34%  '
$ws.Range('E21').Value = '  -0.31%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '80.26'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.76%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.837.37'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0000125'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.93'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.14'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.50%  '
$ws.Range('E29').Value = '  +0.75%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.74'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.04%  '
$ws.Range('E31').Value = '  -1.47%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.57'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.61%  '
$ws.Range('B33').Value = 'Binance-PegBSC-USD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.42%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '26.98'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.99%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.681.23'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.54%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.163'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +1.83%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '8.43'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +2.68%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.23'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.29%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.27'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0907'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.11%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '169.20'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.16%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.941'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '46.94'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.22%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '28.38'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.39%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.72'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.57%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.000280'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.80%  '
$ws.Range('E49').Value = '  -2.38%  '
$ws.Range('E50').Value = '  -1.93%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.86'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.25%  '
